$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '97.094.77'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  +0.49%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.702.15'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  +0.75%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.48'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  -1.37%  '

$ws.Range("E6").Value = '  +2.51%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '660.62'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.423'
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = '  +0.53%  '

$ws.Range("E9").Value = '  -0.03%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.06'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  -2.08%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.698.82'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '  +0.73%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000318'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '  +19.00%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '44.31'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  -3.35%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.209'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  +1.50%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.82'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  +0.81%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.391.56'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  +0.68%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '96.854.14'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  +0.53%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '9.08'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '  +1.99%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.694.80'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  -0.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.04'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  +1.49%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '18.62'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  -0.46%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.504'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  -4.25%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '519.86'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  -1.12%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.43'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  -0.59%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000217'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  +6.67%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.90'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  -2.50%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '102.29'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  +0.03%  '

$ws.Range("E28").Value = '  +16.26%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '13.52'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  +2.96%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.86'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = '  +2.04%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.04'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  -0.15%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  +0.04%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.191'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  +2.42%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.85'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  -1.51%  '

$ws.Range("E35").Value = '  +0.09%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '655.61'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  -0.97%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '32.21'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  -1.16%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.593'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  +0.33%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.85'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  +0.29%  '

$ws.Range("E40").Value = '  +0.01%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.169'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  +4.81%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.85'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '  +5.61%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.05'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  +2.48%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.487'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  +7.39%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '40.14'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  -6.63%  '

$ws.Range("E46").Value = '  +0.41%  '

$ws.Range("E47").Value = '  -1.84%  '

$ws.Range("E48").Value = '  +0.77%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '23.61'
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.73'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  +0.80%  '

$ws.Range("B51").Value = 'OKB'

$ws.Range("C51").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '54.55'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  +0.75%  '
